$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @("29-10-2021", "30-10-2021", "31-10-2021", "01-11-2021", "02-11-2021")

$startRow = 303
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 1).Style = $ws.Cells.Item(2, 1).Style
    $ws.Cells.Item($row, 2).Value = 3623
    $ws.Cells.Item($row, 3).Value = 240
}
